$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the file_name / SEM_no values in the rows where they were entered
# in the wrong columns.
$rows = @(4, 5, 9, 10, 25, 26, 27, 28)
foreach ($r in $rows) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $b
    $ws.Cells.Item($r, 2).Value = $a
}

# Autofit / set explicit column widths for columns A and B.
$ws.Columns.Item(1).ColumnWidth = 19.1640625
$ws.Columns.Item(2).ColumnWidth = 20.1640625

# Freeze the header row and set the view/selection.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F25").Select()

# Adjust the workbook window size.
$excel.ActiveWindow.Width = 26980
$excel.ActiveWindow.Height = 17020
